$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "s"
$ws.Range("C2").Value = "s"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 2
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 34.32
